$d = $word.ActiveDocument

$newText = "Prin această activitate participați în cadrul unei campanii globale de observare și consemnare a celor mai slabe stele vizibile ca metodă de măsurare a poluării luminoase dintr-un anumit loc. Localizând și observând  Constelația Bootes pe cerul nopții și comparând-o cu diagramele stelare, oamenii din întreaga lume vor putea afla în ce măsură iluminatul nocturn din comunitatea lor contribuie la poluarea luminoasă. Contribuțiile dumneavoastră la baza de date online vor facilita o documentare globală privind cerul nocturn observabil."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Prin această activitate*") {
        $r = $p.Range
        # exclude the trailing paragraph mark so we only touch the run text
        $textEnd = $r.End - 1
        $textRange = $d.Range($r.Start, $textEnd)
        $textRange.Delete()

        $insertionPoint = $d.Range($r.Start, $r.Start)
        $insertionPoint.InsertAfter($newText)
        break
    }
}
